$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.512.62'
$ws.Range('E2').Value = '  +0.10%  '
$ws.Range('D3').Value = '1.913.27'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('E4').Value = '  +0.78%  '
$ws.Range('D5').Value = "'325.78"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.18%  '
$ws.Range('E6').Value = '  +0.64%  '
$ws.Range('D7').Value = "'0.4828"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.70%  '
$ws.Range('E8').Value = '  -0.53%  '
$ws.Range('D9').Value = "'0.08174"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.47%  '
$ws.Range('D10').Value = "'1.013"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.33%  '
$ws.Range('E11').Value = '  +4.35%  '
$ws.Range('D12').Value = '1.919.75'
$ws.Range('E12').Value = '  +0.94%  '
$ws.Range('D13').Value = "'6.003"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.29%  '
$ws.Range('D14').Value = "'7.126"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.34%  '
$ws.Range('D15').Value = "'90.47"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.13%  '
$ws.Range('D16').Value = "'0.06805"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.03%  '
$ws.Range('E17').Value = '  +0.73%  '
$ws.Range('E18').Value = '  +0.92%  '
$ws.Range('E19').Value = '  -0.21%  '
$ws.Range('E20').Value = '  +0.61%  '
$ws.Range('D21').Value = '29.534.26'
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('D22').Value = "'5.619"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.50%  '
$ws.Range('D23').Value = "'11.83"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.82%  '
$ws.Range('D24').Value = "'2.182"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.20%  '
$ws.Range('D25').Value = '2.147.60'
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('D26').Value = "'155.68"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.81%  '
$ws.Range('D27').Value = "'6.385"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +8.89%  '
$ws.Range('D28').Value = "'20.06"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.14%  '
$ws.Range('D29').Value = "'2.101"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.50%  '
$ws.Range('D30').Value = "'119.71"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.81%  '
$ws.Range('D31').Value = "'1.030"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.16%  '
$ws.Range('D32').Value = "'0.09570"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.36%  '
$ws.Range('D33').Value = "'5.526"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.45%  '
$ws.Range('D34').Value = "'3.560"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.23%  '
$ws.Range('D35').Value = "'1.395"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.17%  '
$ws.Range('E36').Value = '  +0.32%  '
$ws.Range('D37').Value = "'0.06109"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.15%  '
$ws.Range('D38').Value = "'1.183"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.79%  '
$ws.Range('D39').Value = "'10.83"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.71%  '
$ws.Range('D40').Value = "'0.5943"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.09%  '
$ws.Range('D41').Value = "'7.945"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.53%  '
$ws.Range('E42').Value = '  +0.77%  '
$ws.Range('D43').Value = "'2.460"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.33%  '
$ws.Range('D44').Value = "'1.286"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('D45').Value = "'0.07733"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.29%  '
$ws.Range('D46').Value = "'12.37"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.01%  '
$ws.Range('D47').Value = "'0.5570"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.48%  '
$ws.Range('D48').Value = "'1.945"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.52%  '
$ws.Range('D49').Value = "'115.83"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.36%  '
$ws.Range('E50').Value = '  +1.63%  '
$ws.Range('D51').Value = "'1.054"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.97%  '
